# Apply the "41???? HS circuits data" commit: append new benchmark rows
# (rows 20-34) to the existing sheet, update the view selection, and let
# the used-range / dimension follow automatically from the newly written
# cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> 1-based column index (sheet uses columns B..AD).
$colIndex = @{
    "B" = 2;  "C" = 3;  "D" = 4;  "E" = 5;  "F" = 6;  "G" = 7;  "H" = 8;
    "I" = 9;  "J" = 10; "K" = 11; "L" = 12; "M" = 13; "N" = 14; "O" = 15;
    "P" = 16; "Q" = 17; "R" = 18; "S" = 19; "T" = 20; "U" = 21; "V" = 22;
    "W" = 23; "X" = 24; "Y" = 25; "Z" = 26; "AA" = 27; "AB" = 28; "AC" = 29; "AD" = 30
}

$rowsData = [ordered]@{
    20 = [ordered]@{
        "B" = @{ V = 4110 }
        "C" = @{ V = 10 }
        "D" = @{ V = 56 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "G" = @{ V = 12 }
        "H" = @{ V = 25 }
        "I" = @{ V = 14718566 }
        "J" = @{ V = 33554432 }
        "K" = @{ V = 131300 }
        "L" = @{ V = 1.007598 }
        "M" = @{ V = 0 }
        "N" = @{ V = 0.0000577296039999999; S = 2 }
        "O" = @{ V = 12 }
        "P" = @{ V = 23 }
        "Q" = @{ V = 11636162 }
        "R" = @{ V = 469762048 }
        "S" = @{ V = 818420 }
        "T" = @{ V = 0.998108 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.00000357966400000001; S = 2 }
        "W" = @{ V = 12 }
        "X" = @{ V = 20 }
        "Y" = @{ V = 1027858 }
        "Z" = @{ V = 58720256 }
        "AA" = @{ V = 102760 }
        "AB" = @{ V = 0.998203 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.00000322920900000017; S = 2 }
    }
    21 = [ordered]@{
        "B" = @{ V = 4120 }
        "C" = @{ V = 20 }
        "D" = @{ V = 46 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "G" = @{ V = 12 }
        "H" = @{ V = 25 }
        "I" = @{ V = 24611442 }
        "J" = @{ V = 33554432 }
        "K" = @{ V = 70 }
        "L" = @{ V = 1.312499 }
        "M" = @{ V = 0 }
        "N" = @{ V = 0.097655625001; S = 2 }
        "O" = @{ V = 12 }
        "P" = @{ V = 23 }
        "Q" = @{ V = 18969214 }
        "R" = @{ V = 377487360 }
        "S" = @{ V = 21997 }
        "T" = @{ V = 0.744829 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.065112239241 }
        "W" = @{ V = 12 }
        "X" = @{ V = 22 }
        "Y" = @{ V = 6598495 }
        "Z" = @{ V = 188743680 }
        "AA" = @{ V = 10832 }
        "AB" = @{ V = 0.98109 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.000357588099999999; S = 2 }
    }
    22 = [ordered]@{
        "B" = @{ V = 4130 }
        "C" = @{ V = 30 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "G" = @{ V = 24 }
        "H" = @{ V = 26 }
        "I" = @{ V = 30388334 }
        "J" = @{ V = 67108864 }
        "K" = @{ V = 1 }
        "L" = @{ V = 7.99999 }
        "M" = @{ V = 0 }
        "N" = @{ V = 48.9998600001; S = 2 }
        "O" = @{ V = 24 }
        "P" = @{ V = 24 }
        "Q" = @{ V = 18782562 }
        "R" = @{ V = 721420288 }
        "S" = @{ V = 4241 }
        "T" = @{ V = 0.07029 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.8643606841 }
        "W" = @{ V = 12 }
        "X" = @{ V = 22 }
        "Y" = @{ V = 12956592 }
        "Z" = @{ V = 180355072 }
        "AA" = @{ V = 1074 }
        "AB" = @{ V = 1.012878 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.000165842883999998; S = 2 }
    }
    23 = [ordered]@{
        "B" = @{ V = 4140 }
        "C" = @{ V = 40 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "G" = @{ V = 24 }
        "H" = @{ V = 26 }
        "I" = @{ V = 36919968 }
        "J" = @{ V = 67108864 }
        "K" = @{ V = 0 }
        "L" = @{ V = 0 }
        "M" = @{ V = 0 }
        "N" = @{ V = 1; S = 2 }
        "O" = @{ V = 24 }
        "P" = @{ V = 24 }
        "Q" = @{ V = 24646975 }
        "R" = @{ V = 721420288 }
        "S" = @{ V = 2066 }
        "T" = @{ V = 0.047238 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.907755428644 }
        "W" = @{ V = 24 }
        "X" = @{ V = 23 }
        "Y" = @{ V = 11990607 }
        "Z" = @{ V = 360710144 }
        "AA" = @{ V = 1032 }
        "AB" = @{ V = 1.017534 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.000307441155999997; S = 2 }
    }
    24 = [ordered]@{
        "B" = @{ V = 4150 }
        "C" = @{ V = 50 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "G" = @{ V = 48 }
        "H" = @{ V = 26 }
        "I" = @{ V = 21026295 }
        "J" = @{ V = 67108864 }
        "K" = @{ V = 0 }
        "L" = @{ V = 0 }
        "M" = @{ V = 0 }
        "N" = @{ V = 1; S = 2 }
        "O" = @{ V = 48 }
        "P" = @{ V = 24 }
        "Q" = @{ V = 14077201 }
        "R" = @{ V = 721420288 }
        "S" = @{ V = 1001 }
        "T" = @{ V = 0.023892 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.952786827663999 }
        "W" = @{ V = 48 }
        "X" = @{ V = 24 }
        "Y" = @{ V = 9954197 }
        "Z" = @{ V = 721420288 }
        "AA" = @{ V = 1035 }
        "AB" = @{ V = 1.064321 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.004137191041; S = 2 }
    }
    25 = [ordered]@{
        "B" = @{ V = 4160 }
        "C" = @{ V = 60 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "G" = @{ V = 48 }
        "H" = @{ V = 27 }
        "I" = @{ V = 49055885 }
        "J" = @{ V = 134217728 }
        "K" = @{ V = 0 }
        "L" = @{ V = 0 }
        "M" = @{ V = 0 }
        "N" = @{ V = 1; S = 2 }
        "O" = @{ V = 48 }
        "P" = @{ V = 25 }
        "Q" = @{ V = 33069298 }
        "R" = @{ V = 1442840576 }
        "S" = @{ V = 2023 }
        "T" = @{ V = 0.023051 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.954429348600999 }
        "W" = @{ V = 48 }
        "X" = @{ V = 25 }
        "Y" = @{ V = 27087351 }
        "Z" = @{ V = 1442840576 }
        "AA" = @{ V = 2065 }
        "AB" = @{ V = 1.06054 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.0036650916; S = 2 }
    }
    26 = [ordered]@{
        "B" = @{ V = 4170 }
        "C" = @{ V = 70 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "G" = @{ V = 48 }
        "H" = @{ V = 27 }
        "I" = @{ V = 2510414007 }
        "J" = @{ V = 134217728 }
        "K" = @{ V = 0 }
        "L" = @{ V = 0 }
        "M" = @{ V = 0 }
        "N" = @{ V = 1; S = 2 }
        "O" = @{ V = 48 }
        "P" = @{ V = 27 }
        "Q" = @{ V = 6725957087 }
        "R" = @{ V = 5771362304 }
        "S" = @{ V = 4159 }
        "T" = @{ V = 0.023494 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.953563968036 }
        "W" = @{ V = 48 }
        "X" = @{ V = 27 }
        "Y" = @{ V = 7298131536 }
        "Z" = @{ V = 5771362304 }
        "AA" = @{ V = 4237 }
        "AB" = @{ V = 1.006589 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.0000434149209999994; S = 2 }
    }
    27 = [ordered]@{
        "B" = @{ V = 4180 }
        "C" = @{ V = 80 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "G" = @{ V = 48 }
        "H" = @{ V = 28 }
        "I" = @{ V = 4310433172 }
        "J" = @{ V = 268435456 }
        "K" = @{ V = 0 }
        "L" = @{ V = 0 }
        "M" = @{ V = 0 }
        "N" = @{ V = 1; S = 2 }
        "O" = @{ V = 48 }
        "P" = @{ V = 28 }
        "Q" = @{ V = 235737241 }
        "R" = @{ V = 11542724608 }
        "S" = @{ V = 8117 }
        "T" = @{ V = 0.023242 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.954056190564 }
        "W" = @{ V = 48 }
        "X" = @{ V = 28 }
        "Y" = @{ V = 14947210225 }
        "Z" = @{ V = 11542724608 }
        "AA" = @{ V = 8241 }
        "AB" = @{ V = 1.003171 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.0000100552410000002; S = 2 }
    }
    28 = [ordered]@{
        "B" = @{ V = 4190 }
        "C" = @{ V = 90 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "O" = @{ V = 8 }
        "P" = @{ V = 24 }
        "Q" = @{ V = 56471654 }
        "R" = @{ V = 721420288 }
        "T" = @{ V = 0.022075 }
        "U" = @{ V = 0 }
        "V" = @{ V = 0.956338 }
        "W" = @{ V = 8 }
        "X" = @{ V = 24 }
        "Y" = @{ V = 2720944519 }
        "Z" = @{ V = 721420288 }
        "AB" = @{ V = 0.945309 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.002991; S = 2 }
    }
    29 = [ordered]@{
        "B" = @{ V = 41120 }
        "C" = @{ V = 120 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "W" = @{ V = 48 }
        "X" = @{ V = 27 }
        "Y" = @{ V = 8553062845 }
        "Z" = @{ V = 5771362304 }
        "AA" = @{ V = 2014 }
        "AB" = @{ V = 0.983394 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.000276 }
    }
    30 = [ordered]@{
        "B" = @{ V = 41256 }
        "C" = @{ V = 256 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "W" = @{ V = 48 }
        "X" = @{ V = 28 }
        "Y" = @{ V = 33795072369 }
        "Z" = @{ V = 11542724608 }
        "AA" = @{ V = 8192 }
        "AB" = @{ V = 0.999991 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0 }
    }
    31 = [ordered]@{
        "B" = @{ V = 41512 }
        "C" = @{ V = 512 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "W" = @{ V = 48 }
        "X" = @{ V = 25 }
        "Y" = @{ V = 11693644502 }
        "Z" = @{ V = 1442840576 }
        "AA" = @{ V = 1009 }
        "AB" = @{ V = 0.985334 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.000215 }
    }
    32 = [ordered]@{
        "B" = @{ V = 411024 }
        "C" = @{ V = 1024 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "W" = @{ V = 48 }
        "X" = @{ V = 24 }
        "Y" = @{ V = 20510047491 }
        "Z" = @{ V = 721420288 }
        "AB" = @{ V = 1.033167 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.0011 }
    }
    33 = [ordered]@{
        "B" = @{ V = 412048 }
        "C" = @{ V = 2048 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "W" = @{ V = 48 }
        "X" = @{ V = 24 }
        "Y" = @{ V = 41138387598 }
        "Z" = @{ V = 721420288 }
        "AB" = @{ V = 0.913998 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.007396 }
    }
    34 = [ordered]@{
        "B" = @{ V = 414096 }
        "C" = @{ V = 4096 }
        "D" = @{ V = 43 }
        "E" = @{ V = 1 }
        "F" = @{ V = 0 }
        "W" = @{ V = 48 }
        "X" = @{ V = 25 }
        "Y" = @{ V = 195678067961 }
        "Z" = @{ V = 1442840576 }
        "AA" = @{ V = 1040 }
        "AB" = @{ V = 1.015482 }
        "AC" = @{ V = 0 }
        "AD" = @{ V = 0.00024 }
    }
}

foreach ($r in $rowsData.Keys) {
    $rowData = $rowsData[$r]
    foreach ($col in $rowData.Keys) {
        $cellDef = $rowData[$col]
        $c = $colIndex[$col]
        $target = $ws.Cells.Item([int]$r, $c)
        $target.Value2 = $cellDef.V
        if ($cellDef.ContainsKey("S")) {
            # Reuses the workbook's existing scientific-notation style
            # (numFmtId 11, the same one already used by column N/V/AD).
            $target.NumberFormat = "0.00E+00"
        }
    }
}

# Scroll the view roughly to where the new data was entered and leave the
# active cell / selection where the author left off.
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("V29").Select()
